# Add two new columns, I (I0) and J (IF), to the right of the existing
# H (IP) column on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the header cell formatting (bold, border, centered) from H1 onto
# the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows ----------------------------------------------------------
# Row 2 carries its own, independent values for the new columns.
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 9

# For every remaining data row (3-42), I holds a constant 1 and J simply
# duplicates whatever is already stored in column H on that row.
for ($row = 3; $row -le 42; $row++) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $ws.Cells.Item($row, 8).Value2
}

$excel.ActiveWorkbook.Application.CutCopyMode = 0
